# Apply the diff: swap STAY/CHARTER labels on specific rows of column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "CHARTER"
$ws.Range("A92").Value = "STAY"
$ws.Range("A93").Value = "STAY"
$ws.Range("A113").Value = "CHARTER"
$ws.Range("A114").Value = "CHARTER"
$ws.Range("A141").Value = "CHARTER"
